$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 (shifts existing rows 46-62 down to 47-63,
# and pushes the used range from A1:R62 to A1:R63).
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly price-report entry.
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(46, 3).Value = 'Ñuble'
$ws.Cells.Item(46, 4).Value = 44876
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112026
$ws.Cells.Item(46, 7).Value = 'Haba'
$ws.Cells.Item(46, 8).Value = 'Sin especificar'
$ws.Cells.Item(46, 9).Value = 'Primera'
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 6500
$ws.Cells.Item(46, 12).Value = 7000
$ws.Cells.Item(46, 13).Value = 6750
$ws.Cells.Item(46, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(46, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(46, 16).Value = 270
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = 'Hortaliza'
